$d = $word.ActiveDocument

# --- 1. Remove the "Meta description: ..." paragraph that follows the
#        H1 title paragraph ------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- 2. Replace the final paragraph ("Create a Feature Image Prompt: ...")
#        with two paragraphs: a new bold title paragraph followed by the
#        (now-italic) meta-description text ---------------------------------
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Burning Sun Slot for Free - Unique Features and Mechanics</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Burning Sun slot machine. Play for free with Scatter Pays and Hold the Jackpot Respin feature.</w:t></w:r></w:p>'

$last.Range.InsertXML($xml)

Write-Output "Paragraph count: $($d.Paragraphs.Count)"
